$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Challenges")

# New columns for the Gameweeks import feature.
$ws.Range("S1").Value = "Show Statistics Continuously"
$ws.Range("T1").Value = "Gameweek"

# Sample row values. Leading "'" forces literal text so "true" is stored
# as a string (not auto-coerced to the boolean TRUE), matching the import
# fixture which expects a text column value.
$ws.Range("S2").Value = "'true"
$ws.Range("T2").Value = 1
